$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, pushing the existing row 202 (and
# everything below it) down by one. This mirrors the diff: a brand new
# price record is inserted before the former row 202, and all rows from
# the old 202..265 become 203..266.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new record's data.
$ws.Cells.Item(202, 1).Value = 4
$ws.Cells.Item(202, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(202, 3).Value = "Los Lagos"
$ws.Cells.Item(202, 4).Value = 44588
$ws.Cells.Item(202, 5).Value = 10
$ws.Cells.Item(202, 6).Value = 100112008
$ws.Cells.Item(202, 7).Value = "Coliflor"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 300
$ws.Cells.Item(202, 11).Value = 1500
$ws.Cells.Item(202, 12).Value = 1500
$ws.Cells.Item(202, 13).Value = 1500
$ws.Cells.Item(202, 14).Value = "`$/unidad"
$ws.Cells.Item(202, 15).Value = "Región Metropolitana"
$ws.Cells.Item(202, 16).Value = 1500
$ws.Cells.Item(202, 17).Value = 1
$ws.Cells.Item(202, 18).Value = "Hortaliza"
